# Harvard case classification: recompute the "_old" (Harvard-coded) stats columns
# (Ada_old, Avey_old, Babylon_old, Buoy_old, K health_old, WebMD_old,
#  doctor_MA_old, doctor_NJ_old, doctor_TH_old) and the average_doctor /
# average_doctor_old pair, whose header labels were swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels for the average-doctor columns swap places
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4: stats_for_precision
$ws.Range("E4").Value = 0.44
$ws.Range("F4").Value = 0.064
$ws.Range("G4").Value = 0.253
$ws.Range("N4").Value = 0.429
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.252
$ws.Range("Q4").Value = 0.02
$ws.Range("R4").Value = 0.014
$ws.Range("S4").Value = 0.12
$ws.Range("W4").Value = 0.29
$ws.Range("X4").Value = 0.11
$ws.Range("Y4").Value = 0.332
$ws.Range("AI4").Value = 0.348
$ws.Range("AJ4").Value = 0.08599999999999999
$ws.Range("AK4").Value = 0.292
$ws.Range("AU4").Value = 0.198
$ws.Range("AV4").Value = 0.03
$ws.Range("AW4").Value = 0.174
$ws.Range("BA4").Value = 1.956
$ws.Range("BB4").Value = 0.156
$ws.Range("BC4").Value = 0.395
$ws.Range("BG4").Value = 0.709
$ws.Range("BH4").Value = 0.143
$ws.Range("BI4").Value = 0.378
$ws.Range("BM4").Value = 0.702
$ws.Range("BN4").Value = 0.074
$ws.Range("BO4").Value = 0.273
$ws.Range("BP4").Value = 0.652
$ws.Range("BQ4").Value = 0.699

# Row 5: stats_for_recall
$ws.Range("E5").Value = 0.571
$ws.Range("F5").Value = 0.07199999999999999
$ws.Range("G5").Value = 0.268
$ws.Range("N5").Value = 0.737
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.281
$ws.Range("Q5").Value = 0.01
$ws.Range("R5").Value = 0.003
$ws.Range("S5").Value = 0.051
$ws.Range("W5").Value = 0.275
$ws.Range("X5").Value = 0.107
$ws.Range("Y5").Value = 0.327
$ws.Range("AI5").Value = 0.372
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.308
$ws.Range("AU5").Value = 0.381
$ws.Range("AV5").Value = 0.099
$ws.Range("AW5").Value = 0.315
$ws.Range("BA5").Value = 1.315
$ws.Range("BB5").Value = 0.077
$ws.Range("BC5").Value = 0.277
$ws.Range("BG5").Value = 0.384
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.221
$ws.Range("BM5").Value = 0.551
$ws.Range("BN5").Value = 0.062
$ws.Range("BO5").Value = 0.249
$ws.Range("BP5").Value = 0.438
$ws.Range("BQ5").Value = 0.447

# Row 6: stats_for_f1-score
$ws.Range("E6").Value = 0.497
$ws.Range("N6").Value = 0.542
$ws.Range("Q6").Value = 0.013
$ws.Range("W6").Value = 0.282
$ws.Range("AI6").Value = 0.36
$ws.Range("AU6").Value = 0.261
$ws.Range("BA6").Value = 1.563
$ws.Range("BG6").Value = 0.498
$ws.Range("BM6").Value = 0.617
$ws.Range("BP6").Value = 0.521
$ws.Range("BQ6").Value = 0.542

# Row 7: stats_for_f2-score
$ws.Range("E7").Value = 0.539
$ws.Range("N7").Value = 0.644
$ws.Range("Q7").Value = 0.011
$ws.Range("W7").Value = 0.278
$ws.Range("AI7").Value = 0.367
$ws.Range("AU7").Value = 0.322
$ws.Range("BA7").Value = 1.403
$ws.Range("BG7").Value = 0.423
$ws.Range("BM7").Value = 0.576
$ws.Range("BP7").Value = 0.468
$ws.Range("BQ7").Value = 0.481

# Row 8: stats_for_NDCG
$ws.Range("E8").Value = 0.632
$ws.Range("F8").Value = 0.095
$ws.Range("G8").Value = 0.309
$ws.Range("N8").Value = 0.772
$ws.Range("O8").Value = 0.068
$ws.Range("P8").Value = 0.261
$ws.Range("Q8").Value = 0.01
$ws.Range("S8").Value = 0.077
$ws.Range("W8").Value = 0.318
$ws.Range("X8").Value = 0.123
$ws.Range("Y8").Value = 0.351
$ws.Range("AI8").Value = 0.398
$ws.Range("AJ8").Value = 0.128
$ws.Range("AK8").Value = 0.358
$ws.Range("AU8").Value = 0.325
$ws.Range("AW8").Value = 0.295
$ws.Range("BA8").Value = 1.697
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.352
$ws.Range("BG8").Value = 0.537
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.6860000000000001
$ws.Range("BN8").Value = 0.064
$ws.Range("BO8").Value = 0.254
$ws.Range("BP8").Value = 0.5659999999999999
$ws.Range("BQ8").Value = 0.594

# Row 9: stats_for_M1
$ws.Range("E9").Value = 0.5639999999999999
$ws.Range("F9").Value = 0.246
$ws.Range("G9").Value = 0.496
$ws.Range("N9").Value = 0.667
$ws.Range("O9").Value = 0.222
$ws.Range("P9").Value = 0.471
$ws.Range("W9").Value = 0.218
$ws.Range("X9").Value = 0.17
$ws.Range("Y9").Value = 0.413
$ws.Range("AI9").Value = 0.308
$ws.Range("AJ9").Value = 0.213
$ws.Range("AK9").Value = 0.462
$ws.Range("BA9").Value = 1.615
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.551
$ws.Range("BH9").Value = 0.247
$ws.Range("BI9").Value = 0.497
$ws.Range("BM9").Value = 0.628
$ws.Range("BN9").Value = 0.234
$ws.Range("BO9").Value = 0.483
$ws.Range("BP9").Value = 0.538
$ws.Range("BQ9").Value = 0.5669999999999999

# Row 10: stats_for_M3
$ws.Range("E10").Value = 0.705
$ws.Range("F10").Value = 0.208
$ws.Range("G10").Value = 0.456
$ws.Range("N10").Value = 0.885
$ws.Range("O10").Value = 0.102
$ws.Range("P10").Value = 0.319
$ws.Range("W10").Value = 0.397
$ws.Range("X10").Value = 0.239
$ws.Range("Y10").Value = 0.489
$ws.Range("AI10").Value = 0.436
$ws.Range("AJ10").Value = 0.246
$ws.Range("AK10").Value = 0.496
$ws.Range("AU10").Value = 0.321
$ws.Range("AV10").Value = 0.218
$ws.Range("AW10").Value = 0.467
$ws.Range("BA10").Value = 2.025
$ws.Range("BB10").Value = 0.246
$ws.Range("BC10").Value = 0.496
$ws.Range("BG10").Value = 0.615
$ws.Range("BH10").Value = 0.237
$ws.Range("BI10").Value = 0.487
$ws.Range("BM10").Value = 0.846
$ws.Range("BN10").Value = 0.13
$ws.Range("BO10").Value = 0.361
$ws.Range("BP10").Value = 0.675
$ws.Range("BQ10").Value = 0.715

# Row 11: stats_for_M5
$ws.Range("E11").Value = 0.744
$ws.Range("F11").Value = 0.191
$ws.Range("G11").Value = 0.437
$ws.Range("N11").Value = 0.897
$ws.Range("O11").Value = 0.092
$ws.Range("P11").Value = 0.303
$ws.Range("W11").Value = 0.397
$ws.Range("X11").Value = 0.239
$ws.Range("Y11").Value = 0.489
$ws.Range("AI11").Value = 0.474
$ws.Range("AJ11").Value = 0.249
$ws.Range("AK11").Value = 0.499
$ws.Range("AU11").Value = 0.436
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.025
$ws.Range("BB11").Value = 0.246
$ws.Range("BC11").Value = 0.496
$ws.Range("BG11").Value = 0.615
$ws.Range("BH11").Value = 0.237
$ws.Range("BI11").Value = 0.487
$ws.Range("BM11").Value = 0.846
$ws.Range("BN11").Value = 0.13
$ws.Range("BO11").Value = 0.361
$ws.Range("BP11").Value = 0.675
$ws.Range("BQ11").Value = 0.718

# Row 12: stats_for_position
$ws.Range("E12").Value = 1.448
$ws.Range("F12").Value = 0.834
$ws.Range("G12").Value = 0.913
$ws.Range("N12").Value = 1.423
$ws.Range("O12").Value = 0.751
$ws.Range("P12").Value = 0.867
$ws.Range("W12").Value = 1.613
$ws.Range("X12").Value = 0.5600000000000001
$ws.Range("Y12").Value = 0.748
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.694
$ws.Range("AV12").Value = 2.879
$ws.Range("AW12").Value = 1.697
$ws.Range("BA12").Value = 3.776
$ws.Range("BB12").Value = 0.399
$ws.Range("BC12").Value = 0.631
$ws.Range("BG12").Value = 1.125
$ws.Range("BH12").Value = 0.151
$ws.Range("BI12").Value = 0.389
$ws.Range("BM12").Value = 1.333
$ws.Range("BN12").Value = 0.374
$ws.Range("BO12").Value = 0.611
$ws.Range("BP12").Value = 1.259
$ws.Range("BQ12").Value = 1.281

# Row 13: stats_for_length (x of gs)
$ws.Range("E13").Value = 1.582
$ws.Range("F13").Value = 0.68
$ws.Range("G13").Value = 0.824
$ws.Range("N13").Value = 2.073
$ws.Range("O13").Value = 0.968
$ws.Range("P13").Value = 0.984
$ws.Range("W13").Value = 1.047
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.275
$ws.Range("AV13").Value = 1.042
$ws.Range("AW13").Value = 1.021
$ws.Range("BA13").Value = 2.368
$ws.Range("BB13").Value = 0.297
$ws.Range("BC13").Value = 0.545
$ws.Range("BG13").Value = 0.591
$ws.Range("BH13").Value = 0.077
$ws.Range("BI13").Value = 0.277
$ws.Range("BM13").Value = 0.913
$ws.Range("BN13").Value = 0.294
$ws.Range("BO13").Value = 0.542
$ws.Range("BP13").Value = 0.789
$ws.Range("BQ13").Value = 0.722
